$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("J7").NumberFormat = "General"
